$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a new row at position 2, shifting existing rows 2-11 down to 3-12.
# This preserves the original cell content/types for the shifted rows automatically.
$ws.Rows("2:2").Insert()

# Step 2: Populate the newly inserted row 2 with the new match data (Egyptian Premier).
# Text-like columns (League/Date/Time/Home/Away) must be forced to Text format first
# so Excel does not auto-convert date/time-looking strings into serial numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Egyptian Premier"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2026-01-27"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "12:00:00"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "Al Ahly Cairo"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "Wadi Degla"

# Numeric odds columns for the new row 2
$ws.Range("F2").Value = 1.48
$ws.Range("G2").Value = 1.57
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 10.5
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.2
$ws.Range("O2").Value = 1.38
$ws.Range("P2").Value = 1.74
$ws.Range("Q2").Value = 2.14
$ws.Range("R2").Value = 1.28
$ws.Range("S2").Value = 3.85
$ws.Range("T2").Value = 2.06
$ws.Range("U2").Value = 1.64
$ws.Range("V2").Value = 1.1
$ws.Range("W2").Value = 2.74
$ws.Range("X2").Value = 14
$ws.Range("Y2").Value = 23
$ws.Range("Z2").Value = 80
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 6.8
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 36
$ws.Range("AE2").Value = 210
$ws.Range("AF2").Value = 8
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 34
$ws.Range("AI2").Value = 190
$ws.Range("AJ2").Value = 14.5
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 270
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 360

# Reset row 2 style to Normal so it does not keep the bold/bordered header formatting
# that Excel copies by default when inserting a row above data.
$ws.Range("A2:AO2").Style = "Normal"

# Step 3: Apply the remaining odds adjustments to the shifted rows (now rows 3-12).
$ws.Range("H3").Value = 1.86
$ws.Range("N4").Value = 4.4
$ws.Range("X4").Value = 20
$ws.Range("Z4").Value = 12
$ws.Range("AI4").Value = 65
$ws.Range("F5").Value = 3.35
$ws.Range("G5").Value = 3.4
$ws.Range("H5").Value = 2.24
$ws.Range("I5").Value = 2.26
$ws.Range("K5").Value = 3.9
$ws.Range("N5").Value = 5.4
$ws.Range("P5").Value = 2.48
$ws.Range("Q5").Value = 1.62
$ws.Range("T5").Value = 1.55
$ws.Range("U5").Value = 2.66
$ws.Range("X5").Value = 23
$ws.Range("Y5").Value = 15
$ws.Range("AA5").Value = 30
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 29
$ws.Range("AG5").Value = 15
$ws.Range("AH5").Value = 14.5
$ws.Range("AI5").Value = 29
$ws.Range("AK5").Value = 36
$ws.Range("AM5").Value = 55
$ws.Range("AN5").Value = 24
$ws.Range("F6").Value = 2.18
$ws.Range("I6").Value = 3.9
$ws.Range("F8").Value = 2.56
$ws.Range("G8").Value = 4.1
$ws.Range("H8").Value = 2.64
$ws.Range("I8").Value = 3.25
$ws.Range("J8").Value = 2.68
$ws.Range("K8").Value = 3.65
$ws.Range("P8").Value = 1.38
$ws.Range("F10").Value = 2.04
$ws.Range("G10").Value = 2.2
$ws.Range("F12").Value = 1.93
$ws.Range("G12").Value = 2.32
$ws.Range("J12").Value = 2.82
$ws.Range("P12").Value = 1.57
$ws.Range("Q12").Value = 2.08
